# Natmi following Dr Hou advice:
# Recompute the LR-pairs (Col4a5-Cd93) table with an added "ECs" sending
# cluster and refreshed specificity/weight statistics for every
# sending/target cluster combination (ECs, FAPs, sCs x ECs, FAPs, M2, sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$allRows = @(
    @("ECs", "Col4a5", "Cd93", "ECs", 2, 0.6666666666666666, 0.1821393333333333, 0.5464180000000001, 0.02847959666213861, 0.02847959666213861, 3, 1, 169.629438, 508.888314, 0.7428377317484701, 0.7428377317484702, 30.896192751028, 278.0657347592521, 0.02115571898561434, 0.02115571898561435),
    @("ECs", "Col4a5", "Cd93", "FAPs", 2, 0.6666666666666666, 0.1821393333333333, 0.5464180000000001, 0.02847959666213861, 0.02847959666213861, 3, 1, 0.9848756666666668, 2.954627, 0.004312947180081616, 0.004312947180081616, 0.1793845973428889, 1.614461376086, 0.0001228309961138325, 0.0001228309961138325),
    @("ECs", "Col4a5", "Cd93", "M2", 2, 0.6666666666666666, 0.1821393333333333, 0.5464180000000001, 0.02847959666213861, 0.02847959666213861, 3, 1, 54.620752, 163.862256, 0.2391940691454494, 0.2391940691454494, 9.948587355445335, 89.53728619900801, 0.006812150613238092, 0.006812150613238093),
    @("ECs", "Col4a5", "Cd93", "sCs", 2, 0.6666666666666666, 0.1821393333333333, 0.5464180000000001, 0.02847959666213861, 0.02847959666213861, 3, 1, 3.118221666666667, 9.354665000000001, 0.01365525192599884, 0.01365525192599884, 0.5679508155522223, 5.111557339970001, 0.0003888960671723384, 0.0003888960671723385),
    @("FAPs", "Col4a5", "Cd93", "ECs", 3, 1, 3.954090666666667, 11.862272, 0.618267923195393, 0.6182679231953931, 3, 1, 169.629438, 508.888314, 0.7428377317484701, 0.7428377317484702, 670.7301775877121, 6036.571598289408, 0.4592727416793031, 0.4592727416793033),
    @("FAPs", "Col4a5", "Cd93", "FAPs", 3, 1, 3.954090666666667, 11.862272, 0.618267923195393, 0.6182679231953931, 3, 1, 0.9848756666666668, 2.954627, 0.004312947180081616, 0.004312947180081616, 3.894287681393779, 35.04858913254401, 0.002666556895880487, 0.002666556895880488),
    @("FAPs", "Col4a5", "Cd93", "M2", 3, 1, 3.954090666666667, 11.862272, 0.618267923195393, 0.6182679231953931, 3, 1, 54.620752, 163.862256, 0.2391940691454494, 0.2391940691454494, 215.9754056895147, 1943.778651205632, 0.1478860203712122, 0.1478860203712123),
    @("FAPs", "Col4a5", "Cd93", "sCs", 3, 1, 3.954090666666667, 11.862272, 0.618267923195393, 0.6182679231953931, 3, 1, 3.118221666666667, 9.354665000000001, 0.01365525192599884, 0.01365525192599884, 12.32973118876445, 110.96758069888, 0.008442604248997196, 0.008442604248997196),
    @("sCs", "Col4a5", "Cd93", "ECs", 3, 1, 2.259202333333333, 6.777607, 0.3532524801424683, 0.3532524801424684, 3, 1, 169.629438, 508.888314, 0.7428377317484701, 0.7428377317484702, 383.2272221316219, 3449.044999184598, 0.2624092710835526, 0.2624092710835527),
    @("sCs", "Col4a5", "Cd93", "FAPs", 3, 1, 2.259202333333333, 6.777607, 0.3532524801424683, 0.3532524801424684, 3, 1, 0.9848756666666668, 2.954627, 0.004312947180081616, 0.004312947180081616, 2.225033404176556, 20.025300637589, 0.001523559288087296, 0.001523559288087296),
    @("sCs", "Col4a5", "Cd93", "M2", 3, 1, 2.259202333333333, 6.777607, 0.3532524801424683, 0.3532524801424684, 3, 1, 54.620752, 163.862256, 0.2391940691454494, 0.2391940691454494, 123.3993303668213, 1110.593973301392, 0.08449589816099905, 0.08449589816099906),
    @("sCs", "Col4a5", "Cd93", "sCs", 3, 1, 2.259202333333333, 6.777607, 0.3532524801424683, 0.3532524801424684, 3, 1, 3.118221666666667, 9.354665000000001, 0.01365525192599884, 0.01365525192599884, 7.044693665183889, 63.402242986655, 0.004823751609829308, 0.004823751609829309)
)

for ($r = 0; $r -lt $allRows.Length; $r++) {
    $rowData = $allRows[$r]
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $ws.Cells.Item($r + 2, $c).Value = $rowData[$c - 1]
    }
}